$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2021-10-06"
$ws.Range("B1").Value = "October 2021 (through October 06)"

$ws.Range("L2").Value = 5
$ws.Range("L3").Value = 3
$ws.Range("V3").Value = 1
$ws.Range("AF3").Value = 4
$ws.Range("B4").Value = 4
$ws.Range("BJ4").Value = 2
$ws.Range("B7").Value = 2
$ws.Range("B9").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("B13").Value = 2
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 4
$ws.Range("AF19").Value = 1
$ws.Range("B22").Value = 2
$ws.Range("B28").Value = 1
$ws.Range("B35").Value = 1
$ws.Range("L37").Value = 3
$ws.Range("AZ37").Value = 2
$ws.Range("V39").Value = 1
$ws.Range("B42").Value = 1
$ws.Range("V57").Value = 1
$ws.Range("AF63").Value = 1
